$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# New "Tatsache" (actual hours) figures for the first three tasks:
#   Driver konstante Geschwindigkeit, Obstacles einbauen, Driver Emergency Signal senden
$ws.Range("H2").Value = 35
$ws.Range("H3").Value = 25
$ws.Range("H4").Value = 5

# Move/update the active selection to H3, matching the saved cursor position
$ws.Range("H3").Select()

# Reposition the workbook window (best effort; matches recorded window geometry)
$win = $wb.Windows.Item(1)
$win.Left = -28920
$win.Top = -120
